$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing columns
# (data, widths, styles, data validation) one to the right.
$ws.Columns.Item(1).Insert()

# The inserted column starts out completely empty (no cell records at
# all), so copy the formatting from column B (which now holds what used
# to be column A's styling) into column A before writing any values.
$ws.Range("B1:B8").Copy()
$ws.Range("A1:A8").PasteSpecial(-4122)

# Restore the new column A's width (target width 23). The COM layer's
# ColumnWidth setter adds a fixed +5/6 offset when converting to the
# stored character-width units, so back it out here.
$ws.Columns.Item(1).ColumnWidth = 23 - 0.8333333333333333

# New index column header + values.
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"
$ws.Range("A2").Value = 208
$ws.Range("A3").Value = 209
$ws.Range("A4").Value = 210
$ws.Range("A5").Value = 211
$ws.Range("A6").Value = 212
$ws.Range("A7").Value = 213
$ws.Range("A8").Value = 214

# Re-case / re-word the header row text that changed (cells shifted one
# column right by the insert above, so these are the NEW addresses).
$ws.Range("C1").Value = "REGION"
$ws.Range("Q1").Value = " TARGET COMPLETION DATE "
$ws.Range("R1").Value = "ACTUAL DATE OF COMPLETION"
$ws.Range("S1").Value = "PROJECT ID"
$ws.Range("T1").Value = "CONTRACT ID"
$ws.Range("U1").Value = "ISSUANCE OF INVITATION TO BID"
$ws.Range("V1").Value = "PRE-SUBMISSION CONFERENCE"
$ws.Range("W1").Value = "BID OPENING"
$ws.Range("X1").Value = "ISSUANCE OF RESOLUTION TO AWARD"
$ws.Range("Y1").Value = "ISSUANCE OF NOTICE TO PROCEED"
$ws.Range("Z1").Value = "NAME OF CONTRACTOR"
$ws.Range("AA1").Value = "OTHER REMARKS"
